$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '33.630.56'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.767.51'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '223.20'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.544'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.71'
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.289'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0685'
$ws.Range('E10').Value = '  -3.73%  '
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.020.24'
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.01'
$ws.Range('E13').Value = '  +4.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.776.59'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.669.80'
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('E16').Value = '  -3.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.11'
$ws.Range('E17').Value = '  -2.28%  '
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0774'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.19'
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.52'
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('E23').Value = '  -1.93%  '
$ws.Range('E24').Value = '  -3.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.11'
$ws.Range('E25').Value = '  +1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.05'
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.99'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0510'
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('E32').Value = '  -2.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.48'
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.377.01'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('E37').Value = '  -2.58%  '
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.22'
$ws.Range('E39').Value = '  +5.52%  '
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '77.56'
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.902'
$ws.Range('E42').Value = '  -3.49%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.65'
$ws.Range('E43').Value = '  -2.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.44'
$ws.Range('E44').Value = '  +13.50%  '
$ws.Range('E45').Value = '  +4.08%  '
$ws.Range('E46').Value = '  +13.95%  '
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.61'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.79'
$ws.Range('E49').Value = '  -2.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.920.81'
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('E51').Value = '  +0.10%  '
